$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

# Update the link text in C10 to match the one already used in C9
# (the repository link), replacing the old Netlify site link.
$ws.Range("C10").Value = "https://github.com/contesl/C24172G11"

# Update the active selection to C10 to match the saved view state.
$ws.Range("C10").Select()
